# Revert "typo: Database 2 -> Database 3": change "Database 3" back to
# "Database 2" in the title text of several slides.
#
# Each title shape's first run is edited in place by replacing exactly the
# character span of the existing (old) run text with the new text, which
# keeps the edit confined to a single run and preserves that run's
# formatting (rPr) and all following runs untouched.

$p = $ppt.ActivePresentation

function Set-RunPrefix($slideIndex, $shapeIndex, $oldText, $newText) {
    $tr = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex).TextFrame.TextRange
    $span = $tr.Characters(1, $oldText.Length)
    $span.Text = $newText
}

# Slide 2 - Titel 1
Set-RunPrefix 2 2 "Database 3: ucddb002 128Hz original (+ " "Database 2: ucddb002 128Hz original (+ "

# Slide 3 - Titel 1
Set-RunPrefix 3 2 "Database 3: ucddb002 100Hz (+ " "Database 2: ucddb002 100Hz (+ "

# Slide 4 - Titel 1
Set-RunPrefix 4 2 "Database 3: ucddb025 128Hz original (+ " "Database 2: ucddb025 128Hz original (+ "

# Slide 5 - Titel 1
Set-RunPrefix 5 2 "Database 3: ucddb025 100Hz (+ " "Database 2: ucddb025 100Hz (+ "

# Slide 8 - Titel 1
Set-RunPrefix 8 2 "Database 3: 128Hz - all " "Database 2: 128Hz - all "

# Slide 9 - Titel 1: here the original title is split into two runs,
# "Database 3: " (first run) and "100Hz - all " (second run). The target
# merges them into a single run "Database 2: 100Hz - all " that keeps the
# second run's formatting. Emptying the first run first makes the second
# run become the start of the text range, so writing the combined text into
# it keeps its own rPr instead of inheriting the first run's.
$tr9 = $p.Slides.Item(9).Shapes.Item(2).TextFrame.TextRange
$firstRun = $tr9.Characters(1, "Database 3: ".Length)
$firstRun.Text = ""
$secondRun = $tr9.Characters(1, "100Hz - all ".Length)
$secondRun.Text = "Database 2: 100Hz - all "
